$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 260, shifting the old
# row 260 (Tuna / Primera / 44544) down to row 262.
$ws.Rows.Item(260).Insert()
$ws.Rows.Item(260).Insert()

# Row 259 previously held the "Calameño / Primera / 44544" record; it
# becomes the new "Calameño / Extra / 44595" record.
$ws.Cells.Item(259, 4).Value = 44595
$ws.Cells.Item(259, 9).Value = "Extra"
$ws.Cells.Item(259, 10).Value = 450
$ws.Cells.Item(259, 11).Value = 1500
$ws.Cells.Item(259, 12).Value = 1500
$ws.Cells.Item(259, 13).Value = 1500
$ws.Cells.Item(259, 14).Value = "`$/unidad"
$ws.Cells.Item(259, 16).Value = 1500
$ws.Cells.Item(259, 17).Value = 1

# New row 260: "Tuna / Extra / 44595" record.
$ws.Cells.Item(260, 1).Value = 4
$ws.Cells.Item(260, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(260, 3).Value = "Los Lagos"
$ws.Cells.Item(260, 4).Value = 44595
$ws.Cells.Item(260, 5).Value = 10
$ws.Cells.Item(260, 6).Value = 100112027
$ws.Cells.Item(260, 7).Value = "Melón"
$ws.Cells.Item(260, 8).Value = "Tuna"
$ws.Cells.Item(260, 9).Value = "Extra"
$ws.Cells.Item(260, 10).Value = 450
$ws.Cells.Item(260, 11).Value = 1500
$ws.Cells.Item(260, 12).Value = 1500
$ws.Cells.Item(260, 13).Value = 1500
$ws.Cells.Item(260, 14).Value = "`$/unidad"
$ws.Cells.Item(260, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(260, 16).Value = 1500
$ws.Cells.Item(260, 17).Value = 1
$ws.Cells.Item(260, 18).Value = "Hortaliza"

# New row 261: restore the original "Calameño / Primera / 44544" record
# that used to live at row 259 before it was overwritten above.
$ws.Cells.Item(261, 1).Value = 4
$ws.Cells.Item(261, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(261, 3).Value = "Los Lagos"
$ws.Cells.Item(261, 4).Value = 44544
$ws.Cells.Item(261, 5).Value = 10
$ws.Cells.Item(261, 6).Value = 100112027
$ws.Cells.Item(261, 7).Value = "Melón"
$ws.Cells.Item(261, 8).Value = "Calameño"
$ws.Cells.Item(261, 9).Value = "Primera"
$ws.Cells.Item(261, 10).Value = 200
$ws.Cells.Item(261, 11).Value = 15000
$ws.Cells.Item(261, 12).Value = 15000
$ws.Cells.Item(261, 13).Value = 15000
$ws.Cells.Item(261, 14).Value = "`$/caja 12 unidades"
$ws.Cells.Item(261, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(261, 16).Value = 1250
$ws.Cells.Item(261, 17).Value = 12
$ws.Cells.Item(261, 18).Value = "Hortaliza"

# Row 262 already contains the old row 260 ("Tuna / Primera / 44544")
# record thanks to the row inserts above, so nothing further to do there.
